$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Area"
$ws.Range("B3").Value = "numeric"
$ws.Range("C3").Value = "Area defined by the outer bounday of the fruit measured in the longitudinal section"

$ws.Range("A3").Font.Name = "Calibri"
$ws.Range("A3").Font.Color = 0

$ws.Range("C3").Font.Name = "Calibri"
$ws.Range("C3").Font.Color = 0

$ws.Range("C3").Select()
